$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value = 3966.6667
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 3950
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 3950
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5198
$ws.Range("H65").Value = 3966.6667
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 3950
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 19750
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -25990
$ws.Range("H98").Value = 2000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2000
$ws.Range("N98").Value = -4996
$ws.Range("M98").Value = ""
$ws.Range("H107").Value = 2044.5454
$ws.Range("I107").Value = 2508.3333
$ws.Range("J107").Value = 1488
$ws.Range("K107").Value = 2508.3333
$ws.Range("L107").Value = 1488
$ws.Range("M107").Value = -588.3332999999998
$ws.Range("N107").Value = -5328
$ws.Range("H112").Value = 4804.7236
$ws.Range("J112").Value = 4949.378
$ws.Range("L112").Value = 14848.134
$ws.Range("N112").Value = -17064.134
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("M122").Value = ""
$ws.Range("H125").Value = 1840
$ws.Range("H132").Value = 5297.852
$ws.Range("I132").Value = 5132.304
$ws.Range("K132").Value = 15396.912
$ws.Range("M132").Value = -12866.912
$ws.Range("H138").Value = 2307.8704
$ws.Range("I138").Value = 2636.0625
$ws.Range("J138").Value = 2169.6843
$ws.Range("K138").Value = 7908.1875
$ws.Range("L138").Value = 6509.0529
$ws.Range("M138").Value = -2768.1875
$ws.Range("N138").Value = -16789.0529

$ws = $wb.Worksheets.Item(2)
$ws.Range("I5").Value = 91
$ws.Range("J5").Value = 193
$ws.Range("K5").Value = 91
$ws.Range("L5").Value = 193
$ws.Range("M5").Value = 21
$ws.Range("N5").Value = -417
$ws.Range("H23").Value = 13512.75
$ws.Range("J23").Value = 9746.532999999999
$ws.Range("L23").Value = 9746.532999999999
$ws.Range("N23").Value = -10264.533
$ws.Range("H32").Value = 394186.1
$ws.Range("I32").Value = 445237.56
$ws.Range("J32").Value = 21510.3
$ws.Range("K32").Value = 445237.56
$ws.Range("L32").Value = 21510.3
$ws.Range("M32").Value = -444950.56
$ws.Range("N32").Value = -22084.3
$ws.Range("H37").Value = 12736.363
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 12736.363
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 12736.363
$ws.Range("N37").Value = -13282.363
$ws.Range("M37").Value = ""
$ws.Range("H44").Value = 18398.285
$ws.Range("J44").Value = 18398.285
$ws.Range("L44").Value = 18398.285
$ws.Range("N44").Value = -19374.285
$ws.Range("H110").Value = 2450
$ws.Range("I110").Value = 2450
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2450
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -405
$ws.Range("N110").Value = ""
$ws.Range("H122").Value = 114000
$ws.Range("I122").Value = 202120
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 606360
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -603910
$ws.Range("N122").Value = -16450

$ws = $wb.Worksheets.Item(3)
$ws.Range("I4").Value = 91
$ws.Range("J4").Value = 193
$ws.Range("K4").Value = 91
$ws.Range("L4").Value = 193
$ws.Range("M4").Value = 24
$ws.Range("N4").Value = -423
$ws.Range("H86").Value = 52633584
$ws.Range("I86").Value = 55557588
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 55557588
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -55556465
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 52633584
$ws.Range("I89").Value = 55557588
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 277787940
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -277782324
$ws.Range("N89").Value = -18732

$ws = $wb.Worksheets.Item(4)
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1887
$ws.Range("H6").Value = 1670333.4
$ws.Range("J6").Value = 5500
$ws.Range("L6").Value = 5500
$ws.Range("N6").Value = -5726
$ws.Range("H7").Value = 42.870968
$ws.Range("I7").Value = 213.66667
$ws.Range("J7").Value = 24.571428
$ws.Range("K7").Value = 213.66667
$ws.Range("L7").Value = 24.571428
$ws.Range("M7").Value = -100.66667
$ws.Range("N7").Value = -250.571428
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("H31").Value = 8064.7427
$ws.Range("I31").Value = 2015.8572
$ws.Range("J31").Value = 12097.333
$ws.Range("K31").Value = 2015.8572
$ws.Range("L31").Value = 12097.333
$ws.Range("M31").Value = -1720.8572
$ws.Range("N31").Value = -12687.333
$ws.Range("H34").Value = 8064.7427
$ws.Range("I34").Value = 2015.8572
$ws.Range("J34").Value = 12097.333
$ws.Range("K34").Value = 2015.8572
$ws.Range("L34").Value = 12097.333
$ws.Range("M34").Value = -1813.8572
$ws.Range("N34").Value = -12501.333
$ws.Range("H105").Value = 1885.7142
$ws.Range("I105").Value = 1840
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1840
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -93
$ws.Range("N105").Value = -5494
$ws.Range("H132").Value = 8335599
$ws.Range("I132").Value = 1751.0834
$ws.Range("J132").Value = 20836372
$ws.Range("K132").Value = 5253.2502
$ws.Range("L132").Value = 62509116
$ws.Range("M132").Value = -2723.2502
$ws.Range("N132").Value = -62514176
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item(5)
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -1838
$ws.Range("M25").Value = ""
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1500
$ws.Range("N30").Value = -1704
$ws.Range("M30").Value = ""
$ws.Range("H107").Value = 699.8570999999999
$ws.Range("I107").Value = 733.3333
$ws.Range("J107").Value = 690.7273
$ws.Range("K107").Value = 2199.9999
$ws.Range("L107").Value = 2072.1819
$ws.Range("M107").Value = -279.9998999999998
$ws.Range("N107").Value = -5912.1819

$ws = $wb.Worksheets.Item(6)
$ws.Range("H25").Value = 90009
$ws.Range("J25").Value = 90009
$ws.Range("L25").Value = 90009
$ws.Range("N25").Value = -91067
$ws.Range("H122").Value = 4482.4546
$ws.Range("I122").Value = 4001.1667
$ws.Range("J122").Value = 5060
$ws.Range("K122").Value = 12003.5001
$ws.Range("L122").Value = 15180
$ws.Range("M122").Value = -9553.500100000001
$ws.Range("N122").Value = -20080

$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 881.75
$ws.Range("I46").Value = 556.2
$ws.Range("K46").Value = 556.2
$ws.Range("M46").Value = -368.2
$ws.Range("H128").Value = 39800
$ws.Range("J128").Value = 39800
$ws.Range("L128").Value = 39800
$ws.Range("N128").Value = -49760
$ws.Range("H132").Value = 4936.3125
$ws.Range("I132").Value = 4200.1665
$ws.Range("J132").Value = 5378
$ws.Range("K132").Value = 12600.4995
$ws.Range("L132").Value = 16134
$ws.Range("M132").Value = -10070.4995
$ws.Range("N132").Value = -21194

$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 4340.7036
$ws.Range("I122").Value = 2457.7144
$ws.Range("J122").Value = 4999.75
$ws.Range("K122").Value = 7373.1432
$ws.Range("L122").Value = 14999.25
$ws.Range("M122").Value = -4923.1432
$ws.Range("N122").Value = -19899.25
